$wb = $excel.ActiveWorkbook

# Insert the new "Regiones" sheet as the first sheet in the workbook
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$newSheet.Name = "Regiones"

$newSheet.Range("A1").Value = 'Key'
$newSheet.Range("B1").Value = 'Nombre'
$newSheet.Range("C1").Value = 'Numero'
$newSheet.Range("B2").Value = 'Arica y Parinacota'
$newSheet.Range("A2").Value = 'cl-2740'
$newSheet.Range("C2").Value = 'XV'
$newSheet.Range("B3").Value = 'Tarapacá'
$newSheet.Range("A3").Value = 'cl-ta'
$newSheet.Range("C3").Value = 'I'
$newSheet.Range("B4").Value = 'Antofagasta'
$newSheet.Range("A4").Value = 'cl-an'
$newSheet.Range("C4").Value = 'II'
$newSheet.Range("B5").Value = 'Atacama'
$newSheet.Range("A5").Value = 'cl-at'
$newSheet.Range("C5").Value = 'III'
$newSheet.Range("B6").Value = 'Coquimbo'
$newSheet.Range("A6").Value = 'cl-co'
$newSheet.Range("C6").Value = 'IV'
$newSheet.Range("B7").Value = 'Valparaíso'
$newSheet.Range("A7").Value = 'cl-vs'
$newSheet.Range("C7").Value = 'V'
$newSheet.Range("B8").Value = 'Región Metropolitana'
$newSheet.Range("A8").Value = 'cl-rm'
$newSheet.Range("C8").Value = 'RM'
$newSheet.Range("A9").Value = 'cl-li'
$newSheet.Range("B10").Value = 'Maule'
$newSheet.Range("A10").Value = 'cl-ml'
$newSheet.Range("C10").Value = 'VII'
$newSheet.Range("B11").Value = 'Bío-Bío'
$newSheet.Range("A11").Value = 'cl-bi'
$newSheet.Range("C11").Value = 'VIII'
$newSheet.Range("B12").Value = 'La Araucanía'
$newSheet.Range("A12").Value = 'cl-2730'
$newSheet.Range("C12").Value = 'IX'
$newSheet.Range("B13").Value = 'Los Ríos'
$newSheet.Range("A13").Value = 'cl-ar'
$newSheet.Range("C13").Value = 'XIV'
$newSheet.Range("B14").Value = 'Los Lagos'
$newSheet.Range("A14").Value = 'cl-ll'
$newSheet.Range("C14").Value = 'X'
$newSheet.Range("A15").Value = 'cl-ai'
$newSheet.Range("A16").Value = 'cl-ma'
$newSheet.Range("B9").Value = 'O''Higgins'
$newSheet.Range("C9").Value = 'VI'
$newSheet.Range("B15").Value = 'Aysén'
$newSheet.Range("C15").Value = 'XIV'
$newSheet.Range("B16").Value = 'Magallanes'
$newSheet.Range("C16").Value = 'XII'
$newSheet.Range("B17").Value = 'Ñuble'
$newSheet.Range("C17").Value = '*'

# Autofit the "Nombre" column (column B) now that it is populated, matching
# the width Excel's own best-fit sizing produced for this content
$newSheet.Columns.Item(2).ColumnWidth = 31.3

# Restore the selection on the "Hectareas quemadas" sheet and drop its tab-selected flag
$wsHectareas = $wb.Worksheets.Item("Hectáreas quemadas")
$wsHectareas.Range("B1:Q1").Select()

# Make "Regiones" the active sheet again, with its own selection
$newSheet.Range("C17").Select()
